$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert two new blank rows at positions 679 and 680, shifting the
# existing rows 679-753 down to 681-755 (dimension grows from R753 to R755).
$ws.Rows("679:680").Insert()

function Set-DataRow($RowNum, $Fecha, $Variedad, $Volumen, $PrecioMin, $PrecioMax, $PrecioProm, $Unidad, $Origen, $PrecioKg, $KgUnidades) {
    $ws.Cells.Item($RowNum, 1).Value = 5
    $ws.Cells.Item($RowNum, 2).Value = "Macroferia Regional de Talca"
    $ws.Cells.Item($RowNum, 3).Value = "Maule"
    $ws.Cells.Item($RowNum, 4).Value = $Fecha
    $ws.Cells.Item($RowNum, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $ws.Cells.Item($RowNum, 5).Value = 7
    $ws.Cells.Item($RowNum, 6).Value = 100112043
    $ws.Cells.Item($RowNum, 7).Value = "Pepino ensalada"
    $ws.Cells.Item($RowNum, 8).Value = $Variedad
    $ws.Cells.Item($RowNum, 9).Value = "Primera"
    $ws.Cells.Item($RowNum, 10).Value = $Volumen
    $ws.Cells.Item($RowNum, 11).Value = $PrecioMin
    $ws.Cells.Item($RowNum, 12).Value = $PrecioMax
    $ws.Cells.Item($RowNum, 13).Value = $PrecioProm
    $ws.Cells.Item($RowNum, 14).Value = $Unidad
    $ws.Cells.Item($RowNum, 15).Value = $Origen
    $ws.Cells.Item($RowNum, 16).Value = $PrecioKg
    $ws.Cells.Item($RowNum, 17).Value = $KgUnidades
    $ws.Cells.Item($RowNum, 18).Value = "Hortaliza"
}

# New row 679: Arica y Parinacota lot, $300 volumen, $15000 precios, $250/kg
Set-DataRow 679 45212 "Sin especificar" 300 15000 15000 15000 "`$/caja 60 unidades" "Región de Arica y Parinacota" 250 60

# New row 680: Región del Maule lot, 200 volumen, $17000 precios, $212/kg
Set-DataRow 680 45212 "Sin especificar" 200 17000 17000 17000 "`$/caja 80 unidades" "Región del Maule" 212 80
